# Automatische test-sync: 2025-06-23 18:48:50
# Adds a new "Retour is nog niet verwerkt" log entry to the Logs sheet
# and refreshes the Dashboard summary (Retour / Terugbetaling now has
# 5 occurrences and outranks IT / Technisch probleem which stays at 4).

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Append new row 22 to the Logs sheet ---------------------------------
$row = 22

$logs.Cells.Item($row, 1).Value = "Retour is nog niet verwerkt"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Ik heb iets teruggestuurd maar hoor niks. Wanneer krijg ik mijn geld terug?"
$logs.Cells.Item($row, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item($row, 5).Value = "Beste klant,`nBedankt voor je bericht. Om je vraag over de terugbetaling te kunnen beantwoorden, heb ik wat meer informatie nodig. Zou je mij alsjeblieft de volgende gegevens kunnen verstrekken:`n- Factuurnummer van de geretourneerde aankoop`n- Naam waaronder de bestelling is geplaatst`nMet deze gegevens kan ik je verder helpen en zorgen voor een spoedige afhandeling van de terugbetaling.`nIk kijk uit naar je reactie.`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$logs.Cells.Item($row, 6).Value = "2025-06-23 18:48:14"
$logs.Cells.Item($row, 7).Value = "Ja"

# Avoid leaving a stray custom row-height behind after writing the
# multi-line reply text into column E.
$logs.Rows.Item($row).AutoFit()

# --- Update the Dashboard summary table -----------------------------------
# "Retour / Terugbetaling" now has 5 entries and moves above
# "IT / Technisch probleem" (still 4 entries).
$dashboard.Cells.Item(2, 1).Value = "Retour / Terugbetaling"
$dashboard.Cells.Item(2, 2).Value = 5
$dashboard.Cells.Item(3, 1).Value = "IT / Technisch probleem"
$dashboard.Cells.Item(3, 2).Value = 4

# --- Extend conditional formatting ranges to cover the new row ------------
$logs.Range("D2:D21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D22"))
$logs.Range("G2:G21").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G22"))
